# Generate Report for Handoff
# Updates the localization-status report: the two pending files were
# handed off together (acd3cc98.../ffff5f0b0cd7...) instead of each
# having its own separate handback, so "Latest Target/Handback" info is
# cleared and status flips from "Handed back: in sync with en-US" to
# "Ready for handoff".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---- Overview sheet -------------------------------------------------
$overview.Range("A2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.md"
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

$overview.Range("A3").Value = "ffff5f0b0cd7-2e4d-4fac-8347-038a6aff8eeb.md"
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---- zh-cn sheet ------------------------------------------------------
$zhcn.Range("A2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.md"
$zhcn.Range("B2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.zh-cn.xlf"
$zhcn.Range("D2").Value = "2016-03-09 10:06:00"
$zhcn.Range("E2").Value = ""
$zhcn.Range("F2").Value = ""
$zhcn.Range("G2").Value = "0001-01-01 00:00:00"

$zhcn.Range("A3").Value = "ffff5f0b0cd7-2e4d-4fac-8347-038a6aff8eeb.md"
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-09 10:06:00"
$zhcn.Range("E3").Value = ""
$zhcn.Range("F3").Value = ""
$zhcn.Range("G3").Value = "0001-01-01 00:00:00"

# ---- de-de sheet ------------------------------------------------------
$dede.Range("A2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.md"
$dede.Range("B2").Value = "Ready for handoff"
$dede.Range("C2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.de-de.xlf"
$dede.Range("D2").Value = "2016-03-09 10:06:06"
$dede.Range("E2").Value = ""
$dede.Range("F2").Value = ""
$dede.Range("G2").Value = "0001-01-01 00:00:00"

$dede.Range("A3").Value = "ffff5f0b0cd7-2e4d-4fac-8347-038a6aff8eeb.md"
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.de-de.xlf"
$dede.Range("D3").Value = "2016-03-09 10:06:06"
$dede.Range("E3").Value = ""
$dede.Range("F3").Value = ""
$dede.Range("G3").Value = "0001-01-01 00:00:00"
